$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2399.6667
$ws.Range("J40").Value = 2499.6
$ws.Range("L40").Value = 2499.6
$ws.Range("N40").Value = -2849.6
$ws.Range("H69").Value = 8462.25
$ws.Range("J69").Value = 8462.25
$ws.Range("L69").Value = 25386.75
$ws.Range("N69").Value = -27134.75
$ws.Range("H72").Value = 8462.25
$ws.Range("J72").Value = 8462.25
$ws.Range("L72").Value = 76160.25
$ws.Range("N72").Value = -84896.25
$ws.Range("H92").Value = 45829.727
$ws.Range("I92").Value = 254.28572
$ws.Range("K92").Value = 254.28572
$ws.Range("M92").Value = 993.71428
$ws.Range("H101").Value = 634.75
$ws.Range("I101").Value = 429.66666
$ws.Range("J101").Value = 1250
$ws.Range("K101").Value = 1288.99998
$ws.Range("L101").Value = 3750
$ws.Range("M101").Value = 333.0000199999999
$ws.Range("N101").Value = -6994
$ws.Range("H137").Value = 2397.4783
$ws.Range("I137").Value = 1542.4615
$ws.Range("K137").Value = 4627.3845
$ws.Range("M137").Value = -2077.3845
$ws.Range("H141").Value = 5796.25
$ws.Range("J141").Value = 5790
$ws.Range("L141").Value = 17370
$ws.Range("N141").Value = -27730

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 111113890
$ws.Range("I61").Value = 166668080
$ws.Range("K61").Value = 166668080
$ws.Range("M61").Value = -166667868
$ws.Range("H63").Value = 150012480
$ws.Range("I63").Value = 250002130
$ws.Range("K63").Value = 250002130
$ws.Range("M63").Value = -250001444
$ws.Range("H66").Value = 150012480
$ws.Range("I66").Value = 250002130
$ws.Range("K66").Value = 1250010650
$ws.Range("M66").Value = -1250007218
$ws.Range("H74").Value = 2380
$ws.Range("I74").Value = 2299.75
$ws.Range("K74").Value = 2299.75
$ws.Range("M74").Value = -1425.75
$ws.Range("H77").Value = 2380
$ws.Range("I77").Value = 2299.75
$ws.Range("K77").Value = 11498.75
$ws.Range("M77").Value = -7130.75
$ws.Range("H136").Value = 111113890
$ws.Range("I136").Value = 166668080
$ws.Range("K136").Value = 500004240
$ws.Range("M136").Value = -500001690

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 29999
$ws.Range("I75").Value = 29998.5
$ws.Range("K75").Value = 29998.5
$ws.Range("M75").Value = -29062.5
$ws.Range("H78").Value = 29999
$ws.Range("I78").Value = 29998.5
$ws.Range("K78").Value = 89995.5
$ws.Range("M78").Value = -85315.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2778.7642
$ws.Range("I31").Value = 1865.9429
$ws.Range("J31").Value = 3370.4075
$ws.Range("K31").Value = 1865.9429
$ws.Range("L31").Value = 3370.4075
$ws.Range("M31").Value = -1570.9429
$ws.Range("N31").Value = -3960.4075
$ws.Range("H34").Value = 2778.7642
$ws.Range("I34").Value = 1865.9429
$ws.Range("J34").Value = 3370.4075
$ws.Range("K34").Value = 1865.9429
$ws.Range("L34").Value = 3370.4075
$ws.Range("M34").Value = -1663.9429
$ws.Range("N34").Value = -3774.4075
$ws.Range("H58").Value = 2022
$ws.Range("J58").Value = 2583.3333
$ws.Range("L58").Value = 2583.3333
$ws.Range("N58").Value = -2989.3333
$ws.Range("H86").Value = 5136
$ws.Range("I86").Value = 3942.3333
$ws.Range("K86").Value = 3942.3333
$ws.Range("M86").Value = -2819.3333
$ws.Range("H89").Value = 5136
$ws.Range("I89").Value = 3942.3333
$ws.Range("K89").Value = 19711.6665
$ws.Range("M89").Value = -14095.6665
$ws.Range("H105").Value = 898.1429000000001
$ws.Range("I105").Value = 741.1667
$ws.Range("K105").Value = 741.1667
$ws.Range("M105").Value = 1005.8333
$ws.Range("H122").Value = 2296
$ws.Range("J122").Value = 2250
$ws.Range("L122").Value = 6750
$ws.Range("N122").Value = -11650
$ws.Range("H134").Value = 2610
$ws.Range("I134").Value = 2441.7778
$ws.Range("K134").Value = 7325.3334
$ws.Range("M134").Value = -4790.3334
$ws.Range("H136").Value = 2022
$ws.Range("J136").Value = 2583.3333
$ws.Range("L136").Value = 7749.999899999999
$ws.Range("N136").Value = -12849.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 512.17645
$ws.Range("J122").Value = 480.69232
$ws.Range("L122").Value = 4326.23088
$ws.Range("N122").Value = -9226.230879999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 650336.7
$ws.Range("J4").Value = 475505
$ws.Range("L4").Value = 475505
$ws.Range("N4").Value = -475731
$ws.Range("H16").Value = 2537.1
$ws.Range("I16").Value = 2266.5
$ws.Range("J16").Value = 2943
$ws.Range("K16").Value = 2266.5
$ws.Range("L16").Value = 2943
$ws.Range("M16").Value = -2096.5
$ws.Range("N16").Value = -3283
$ws.Range("H22").Value = 728.4286
$ws.Range("I22").Value = 830.6667
$ws.Range("J22").Value = 651.75
$ws.Range("K22").Value = 830.6667
$ws.Range("L22").Value = 651.75
$ws.Range("M22").Value = -535.6667
$ws.Range("N22").Value = -1241.75
$ws.Range("H27").Value = 728.4286
$ws.Range("I27").Value = 830.6667
$ws.Range("J27").Value = 651.75
$ws.Range("K27").Value = 830.6667
$ws.Range("L27").Value = 651.75
$ws.Range("M27").Value = -723.6667
$ws.Range("N27").Value = -865.75
$ws.Range("H28").Value = 650336.7
$ws.Range("J28").Value = 475505
$ws.Range("L28").Value = 475505
$ws.Range("N28").Value = -475969
$ws.Range("H37").Value = 650336.7
$ws.Range("J37").Value = 475505
$ws.Range("L37").Value = 475505
$ws.Range("N37").Value = -475719
$ws.Range("H46").Value = 973
$ws.Range("I46").Value = 973
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 973
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -785
$ws.Range("N46").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 36598.8
$ws.Range("I64").Value = 29998.5
$ws.Range("K64").Value = 29998.5
$ws.Range("M64").Value = -29750.5
$ws.Range("H67").Value = 36598.8
$ws.Range("I67").Value = 29998.5
$ws.Range("K67").Value = 29998.5
$ws.Range("M67").Value = -29140.5
$ws.Range("H128").Value = 54999
$ws.Range("J128").Value = 54999
$ws.Range("L128").Value = 54999
$ws.Range("N128").Value = -64959
$ws.Range("H132").Value = 7155.5713
$ws.Range("I132").Value = 6925.364
$ws.Range("K132").Value = 20776.092
$ws.Range("M132").Value = -18246.092
